# Better detection of numerical data types. Not only System.Double.
# Insert three new rows ("Decimal Number:", "Float Number:", "Double Number:")
# right after the "TimeSpan:" row (row 7), keeping the existing blank-row
# separator pattern used throughout the sheet (row 8 blank, new rows 9-11,
# row 12 blank, then the previously-existing content continues at row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 9 downward by 4 rows (3 new data rows + 1 blank
# separator row), matching the gap pattern already used between sections.
$ws.Range("A9:A12").EntireRow.Insert()

$ws.Range("B9").Value = "Decimal Number:"
$ws.Range("C9").Value = 123.45

$ws.Range("B10").Value = "Float Number:"
$ws.Range("C10").Value = 123.45

$ws.Range("B11").Value = "Double Number:"
$ws.Range("C11").Value = 123.45
